$d = $word.ActiveDocument

# The paragraph currently reads:
#   "Jhbuogouygi7tftttttttttttttttf0i9upi9"
# The edit inserts "sfdsfs " right before the trailing "i9upi9", and the
# existing _GoBack bookmark ends up sitting between the two pieces of text
# (i.e. the "i9upi9" tail moves into its own run after the bookmark).

$insertText = "sfdsfs "
$tailText = "i9upi9"
$bookmarkName = "_GoBack"

# Locate the "i9upi9" tail that the new text is inserted in front of.
$tailRange = $d.Content
$found = $tailRange.Find.Execute($tailText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPos = $tailRange.Start

# Type the new text immediately before that tail.
$insertionPoint = $d.Range($splitPos, $splitPos)
$insertionPoint.InsertBefore($insertText)

# The split point between the two runs is now shifted by the inserted text.
$newSplitPos = $splitPos + $insertText.Length

# Re-seat the _GoBack bookmark exactly at that split point so it ends up
# between the two runs in the saved XML.
if ($d.Bookmarks.Exists($bookmarkName)) {
    $d.Bookmarks($bookmarkName).Delete()
}
$bookmarkRange = $d.Range($newSplitPos, $newSplitPos)
$d.Bookmarks.Add($bookmarkName, $bookmarkRange)

Write-Output ("Final text: [" + $d.Content.Text + "]")
